$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.051938
$ws.Range("H2").Value = 0.155814
$ws.Range("I2").Value = 0.1172837182974765
$ws.Range("J2").Value = 0.1172837182974765
$ws.Range("M2").Value = 34.32490066666667
$ws.Range("N2").Value = 102.974702
$ws.Range("O2").Value = 0.7676290729251072
$ws.Range("P2").Value = 0.7676290729251072
$ws.Range("Q2").Value = 1.782766690825334
$ws.Range("R2").Value = 16.044900217428
$ws.Range("S2").Value = 0.09003039194590128
$ws.Range("T2").Value = 0.09003039194590128
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.051938
$ws.Range("H3").Value = 0.155814
$ws.Range("I3").Value = 0.1172837182974765
$ws.Range("J3").Value = 0.1172837182974765
$ws.Range("M3").Value = 9.765320666666668
$ws.Range("O3").Value = 0.2183879313436534
$ws.Range("P3").Value = 0.2183879313436534
$ws.Range("Q3").Value = 0.5071912247853334
$ws.Range("R3").Value = 4.564721023068
$ws.Range("S3").Value = 0.02561334861927767
$ws.Range("T3").Value = 0.02561334861927768
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.051938
$ws.Range("H4").Value = 0.155814
$ws.Range("I4").Value = 0.1172837182974765
$ws.Range("J4").Value = 0.1172837182974765
$ws.Range("O4").Value = 0.01398299573123946
$ws.Range("P4").Value = 0.01398299573123946
$ws.Range("Q4").Value = 0.03247456344066667
$ws.Range("R4").Value = 0.292271070966
$ws.Range("S4").Value = 0.001639977732297504
$ws.Range("T4").Value = 0.001639977732297504
$ws.Range("I5").Value = 0.5474768201053503
$ws.Range("J5").Value = 0.5474768201053503
$ws.Range("M5").Value = 34.32490066666667
$ws.Range("N5").Value = 102.974702
$ws.Range("O5").Value = 0.7676290729251072
$ws.Range("P5").Value = 0.7676290729251072
$ws.Range("Q5").Value = 8.321900542130003
$ws.Range("R5").Value = 74.89710487917002
$ws.Range("S5").Value = 0.4202591238654557
$ws.Range("T5").Value = 0.4202591238654557
$ws.Range("I6").Value = 0.5474768201053503
$ws.Range("J6").Value = 0.5474768201053503
$ws.Range("M6").Value = 9.765320666666668
$ws.Range("O6").Value = 0.2183879313436534
$ws.Range("P6").Value = 0.2183879313436534
$ws.Range("R6").Value = 21.30797852127001
$ws.Range("S6").Value = 0.1195623302014089
$ws.Range("T6").Value = 0.1195623302014089
$ws.Range("I7").Value = 0.5474768201053503
$ws.Range("J7").Value = 0.5474768201053503
$ws.Range("O7").Value = 0.01398299573123946
$ws.Range("P7").Value = 0.01398299573123946
$ws.Range("S7").Value = 0.007655366038485665
$ws.Range("T7").Value = 0.007655366038485665
$ws.Range("I8").Value = 0.3352394615971734
$ws.Range("J8").Value = 0.3352394615971734
$ws.Range("M8").Value = 34.32490066666667
$ws.Range("N8").Value = 102.974702
$ws.Range("O8").Value = 0.7676290729251072
$ws.Range("P8").Value = 0.7676290729251072
$ws.Range("Q8").Value = 5.095794661538446
$ws.Range("R8").Value = 45.862151953846
$ws.Range("S8").Value = 0.2573395571137503
$ws.Range("T8").Value = 0.2573395571137503
$ws.Range("I9").Value = 0.3352394615971734
$ws.Range("J9").Value = 0.3352394615971734
$ws.Range("M9").Value = 9.765320666666668
$ws.Range("O9").Value = 0.2183879313436534
$ws.Range("P9").Value = 0.2183879313436534
$ws.Range("S9").Value = 0.07321225252296684
$ws.Range("T9").Value = 0.07321225252296686
$ws.Range("I10").Value = 0.3352394615971734
$ws.Range("J10").Value = 0.3352394615971734
$ws.Range("O10").Value = 0.01398299573123946
$ws.Range("P10").Value = 0.01398299573123946
$ws.Range("Q10").Value = 0.09282409631522223
$ws.Range("R10").Value = 0.835416866837
$ws.Range("S10").Value = 0.004687651960456289
$ws.Range("T10").Value = 0.00468765196045629
